# Auto-generated script applying 2024-10-18 YTD violent crime data updates
# to the Chicago violent-crime-ytd workbook.
#
# For each affected worksheet, the relevant year-to-date cell(s) in column
# K (2024) -- and, for a few sheets, small corrections to prior-year totals --
# are updated to match the newly published data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6425
$ws.Range("K3").Value = 6638
$ws.Range("F4").Value = 1581
$ws.Range("I4").Value = 1487
$ws.Range("K4").Value = 1377
$ws.Range("K5").Value = 473
$ws.Range("J6").Value = 8500
$ws.Range("K6").Value = 7302
$ws.Range("F7").Value = 19403
$ws.Range("I7").Value = 20785
$ws.Range("J7").Value = 23207
$ws.Range("K7").Value = 22215

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 160
$ws.Range("K7").Value = 663
$ws.Range("K8").Value = 1459
$ws.Range("K10").Value = 131
$ws.Range("K11").Value = 413
$ws.Range("K15").Value = 229
$ws.Range("K18").Value = 146
$ws.Range("K19").Value = 648
$ws.Range("K20").Value = 530
$ws.Range("K23").Value = 220
$ws.Range("K24").Value = 70
$ws.Range("K27").Value = 208
$ws.Range("K29").Value = 1199
$ws.Range("K31").Value = 248
$ws.Range("K33").Value = 973
$ws.Range("K36").Value = 283
$ws.Range("K40").Value = 50
$ws.Range("K41").Value = 157
$ws.Range("K42").Value = 821
$ws.Range("J43").Value = 198
$ws.Range("K47").Value = 150
$ws.Range("K48").Value = 277
$ws.Range("K50").Value = 105
$ws.Range("F63").Value = 160
$ws.Range("I63").Value = 187
$ws.Range("K63").Value = 64
$ws.Range("K67").Value = 870
$ws.Range("K71").Value = 67
$ws.Range("K75").Value = 70
$ws.Range("K90").Value = 207
$ws.Range("K94").Value = 296
$ws.Range("K95").Value = 364
$ws.Range("K96").Value = 236
$ws.Range("K98").Value = 112
$ws.Range("K99").Value = 368
$ws.Range("F101").Value = 19403
$ws.Range("I101").Value = 20785
$ws.Range("J101").Value = 23207
$ws.Range("K101").Value = 22215

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 45
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 215
$ws.Range("K7").Value = 663

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 107
$ws.Range("K7").Value = 413

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 400
$ws.Range("K3").Value = 444
$ws.Range("K4").Value = 82
$ws.Range("K7").Value = 1459

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 248
$ws.Range("K5").Value = 27
$ws.Range("K6").Value = 300
$ws.Range("K7").Value = 973

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 364

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 175
$ws.Range("K6").Value = 186

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 152
$ws.Range("K7").Value = 368

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 246
$ws.Range("K7").Value = 870

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 70
$ws.Range("K6").Value = 235

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 339
$ws.Range("K3").Value = 429
$ws.Range("K6").Value = 345
$ws.Range("K7").Value = 1199

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 130
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 213
$ws.Range("K7").Value = 648

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 223
$ws.Range("K3").Value = 250
$ws.Range("K7").Value = 821

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 183
$ws.Range("K3").Value = 172
$ws.Range("K4").Value = 24
$ws.Range("K7").Value = 530

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 110
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 283

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 77
$ws.Range("K7").Value = 296

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 22
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 207

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 50
